# Update the "view/interest count" figures in column F on both the
# "展览" (Exhibitions) and "全部类型" (All types) sheets, reflecting a
# refreshed data scrape (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 4981
$ws1.Range("F6").Value  = 4981
$ws1.Range("F13").Value = 8175
$ws1.Range("F14").Value = 8175
$ws1.Range("F18").Value = 592
$ws1.Range("F19").Value = 2492
$ws1.Range("F21").Value = 2288
$ws1.Range("F33").Value = 6776
$ws1.Range("F45").Value = 68

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 4981
$ws4.Range("F5").Value  = 4981
$ws4.Range("F12").Value = 8175
$ws4.Range("F13").Value = 8175
$ws4.Range("F16").Value = 592
$ws4.Range("F17").Value = 2492
$ws4.Range("F21").Value = 2288
$ws4.Range("F34").Value = 6776
$ws4.Range("F44").Value = 68
